$wb = $excel.ActiveWorkbook

# --- Sheet "Logs": append a new row with the latest processed e-mail ---
$logs = $wb.Worksheets.Item("Logs")

$logs.Range("A38").Value = "Offerte voor zakelijke samenwerking"
$logs.Range("B38").Value = "mailmind.test@zohomail.eu"
$logs.Range("C38").Value = "Kunt u mij een offerte sturen voor 100 stuks product X?"
$logs.Range("D38").Value = "Offerte / Prijsaanvraag"
$logs.Range("F38").Value = "2025-06-19 22:11:13"
$logs.Range("G38").Value = "Nee"

# Extend the conditional-formatting ranges so the new row is covered too.
$catFormats = $logs.Range("D2:D37").FormatConditions
for ($i = 1; $i -le $catFormats.Count; $i++) {
    $catFormats.Item($i).ModifyAppliesToRange($logs.Range("D2:D38"))
}

$answeredFormats = $logs.Range("G2:G37").FormatConditions
for ($i = 1; $i -le $answeredFormats.Count; $i++) {
    $answeredFormats.Item($i).ModifyAppliesToRange($logs.Range("G2:G38"))
}

# --- Sheet "Dashboard": the category counts table swapped two rows ---
$dash = $wb.Worksheets.Item("Dashboard")

$dash.Range("A5").Value = "Offerte / Prijsaanvraag"
$dash.Range("B5").Value = 4

$dash.Range("A6").Value = "Klacht / Probleem"
$dash.Range("B6").Value = 4
